$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Re-style row 20 (A20:B20) from the old "blank wrap / no border" look to
#    the normal bordered data-row look used throughout the sheet, by copying
#    the format from an existing data row (A2:B2) before writing values.
# ---------------------------------------------------------------------------
$ws.Range("A2:B2").Copy()
$ws.Range("A20:B20").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Build out the new data rows (19 partially, 21-23, 26-30) using the same
#    bordered/wrapped data style, copied from an existing data row so the
#    engine reuses the existing cell-style entry instead of fabricating new
#    (duplicate) ones.
# ---------------------------------------------------------------------------
$ws.Range("A2:B2").Copy()
$ws.Range("A21:B23").PasteSpecial(-4122)
$ws.Range("A26:B30").PasteSpecial(-4122)

# A19 already carries the data style (only B19 had a value before) - leave it.

# ---------------------------------------------------------------------------
# 3) New section header row 25 ("city_selection_screen.dart"): merge first,
#    then paste the header format (bold font, yellow fill, border) from the
#    existing A1:B1 header so merging doesn't split the border, then write
#    the text.
# ---------------------------------------------------------------------------
$ws.Range("A25:B25").Merge()
$ws.Range("A1:B1").Copy()
$ws.Range("A25:B25").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Write cell values. Order matters: it reproduces the exact order in
#    which new shared strings were appended by the original author.
# ---------------------------------------------------------------------------
$ws.Range("A19").Value = "No photos to display."
$ws.Range("B19").Value = "Không có dữ liệu ảnh để hiển thị."

$ws.Range("A20").Value = "Add item"
$ws.Range("B20").Value = "Thêm món đồ"

$ws.Range("A21").Value = "Previous"
$ws.Range("A22").Value = "Next"
$ws.Range("B21").Value = "Trước"
$ws.Range("B22").Value = "Sau"
$ws.Range("B23").Value = "Lưu tất cả"
$ws.Range("A23").Value = "Save all"

$ws.Range("A25").Value = "city_selection_screen.dart"

$ws.Range("B26").Value = "Chọn vị trí"
$ws.Range("A26").Value = "Select location"

$ws.Range("A27").Value = "Save"
$ws.Range("B27").Value = "Lưu"

$ws.Range("A28").Value = "Auto-detect"
$ws.Range("A29").Value = "Manually"
$ws.Range("B29").Value = "Chọn thủ công"
$ws.Range("B28").Value = "Tự động phát hiện"

$ws.Range("A30").Value = "Search city/location…"
$ws.Range("B30").Value = "Tìm kiếm địa danh…"

# ---------------------------------------------------------------------------
# 5) Row heights: every data row in this sheet uses 15.75pt - make sure the
#    freshly-touched rows match. (Row 24 is intentionally left blank/unused,
#    so it is skipped to avoid materialising a phantom empty row.)
# ---------------------------------------------------------------------------
$ws.Range("A19:B23").RowHeight = 15.75
$ws.Range("A25:B30").RowHeight = 15.75

# ---------------------------------------------------------------------------
# 6) Selection / scroll position, mirroring the diff's sheetView change.
# ---------------------------------------------------------------------------
$ws.Range("A31").Select()
